$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every touched cell keeps its original text formatting (no auto number/date conversion)
$updates = @{
    'D2' = '30.831.17'
    'E2' = '  +0.47%  '
    'D3' = '1.917.81'
    'E3' = '  +1.56%  '
    'E4' = '  +0.23%  '
    'D5' = '241.33'
    'E5' = '  -2.35%  '
    'E6' = '  +0.24%  '
    'D7' = '0.4910'
    'E7' = '  -0.45%  '
    'D8' = '0.2971'
    'E8' = '  +0.51%  '
    'D9' = '0.06761'
    'E9' = '  -0.65%  '
    'D10' = '1.889.56'
    'E10' = '  +0.03%  '
    'E11' = '  -0.67%  '
    'D12' = '0.07318'
    'E13' = '  +2.16%  '
    'D14' = '89.41'
    'E14' = '  -2.19%  '
    'D15' = '0.6714'
    'E15' = '  -1.00%  '
    'D16' = '30.806.15'
    'E16' = '  +0.43%  '
    'D17' = '0.000007981'
    'E17' = '  -0.10%  '
    'D18' = '13.56'
    'E18' = '  +2.46%  '
    'E19' = '  +0.28%  '
    'D20' = '2.154.01'
    'E20' = '  +1.02%  '
    'D21' = '1.005'
    'E21' = '  +0.37%  '
    'D22' = '5.245'
    'E22' = '  +8.61%  '
    'D23' = '202.92'
    'E23' = '  +9.18%  '
    'D24' = '6.270'
    'E24' = '  +3.63%  '
    'D25' = '9.651'
    'E25' = '  +3.27%  '
    'D26' = '160.72'
    'E26' = '  +3.01%  '
    'D27' = '18.91'
    'E27' = '  -1.02%  '
    'D28' = '1.975'
    'E28' = '  +3.67%  '
    'D29' = '1.431'
    'E29' = '  +2.24%  '
    'D30' = '4.351'
    'E30' = '  +1.05%  '
    'D31' = '0.09208'
    'E31' = '  +2.37%  '
    'B32' = 'Hedera'
    'C32' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D32' = '0.05549'
    'E32' = '  +6.86%  '
    'B33' = 'Filecoin'
    'C33' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D33' = '4.071'
    'E33' = '  +1.59%  '
    'D34' = '0.7490'
    'E34' = '  +0.63%  '
    'D35' = '1.121'
    'E35' = '  +0.70%  '
    'D36' = '2.719'
    'E36' = '  -0.69%  '
    'E37' = '  +1.24%  '
    'D38' = '2.729'
    'E38' = '  +2.19%  '
    'D39' = '0.9258'
    'E39' = '  -1.42%  '
    'D40' = '2.079'
    'E40' = '  -3.67%  '
    'D41' = '0.4498'
    'E41' = '  +1.74%  '
    'D42' = '72.86'
    'E42' = '  +25.88%  '
    'D43' = '107.56'
    'E43' = '  +2.07%  '
    'D44' = '5.926'
    'E44' = '  +2.62%  '
    'E45' = '  +0.91%  '
    'D46' = '0.1394'
    'E46' = '  +4.02%  '
    'D47' = '7.709'
    'E47' = '  +1.24%  '
    'D48' = '36.53'
    'E48' = '  +8.99%  '
    'B49' = 'Cronos'
    'C49' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D49' = '0.06019'
    'E49' = '  +2.95%  '
    'B50' = 'EnergySwap'
    'C50' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D50' = '9.033'
    'E50' = '  +3.72%  '
    'D51' = '0.4059'
    'E51' = '  +3.23%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}

$wb.Save()
